{"js": "// Word JS API (Office.js) script \u2014 body of `async (context) => { ... }`.\n//\n// Target edit (per the commit's XML diff):\n//  1. \"...suscrito con fecha  11 de marzo de 2022 entre...\"  -> \"...21 de marzo de 2022...\"\n//  2. The table cell holding \"Almod\u00f3var del Campo\" (LOCALIDAD DE RESIDENCIA,\n//     row for \"D\u00edez Vi\u00f1as Malena\") is emptied out.\n//  3. The entire table row for \"Jim\u00e9nez Coello Daniel\" is removed.\n//  4. \"En Puertollano a  11  de marzo  2022\" -> \"En Puertollano a  21  de marzo  2022\"\n\nconst body = context.document.body;\n\n// 1) Fix the date in the intro paragraph (unique match incl. surrounding text).\nconst introMatches = body.search(\"fecha  11 de marzo de 2022\", { matchCase: true });\nintroMatches.load(\"items\");\nawait context.sync();\nif (introMatches.items.length > 0) {\n  introMatches.items[0].insertText(\"fecha  21 de marzo de 2022\", \"Replace\");\n}\n\n// 2) & 3) Work on the student-roster table (2nd table in the document body).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst rosterTable = tables.items[1];\n\n// 2) Empty the \"Almod\u00f3var del Campo\" cell \u2014 row index 2 (\"D\u00edez Vi\u00f1as Malena\"),\n//    column index 2 (\"LOCALIDAD DE RESIDENCIA DEL ALUMNO/A (**)\").\nconst localidadCell = rosterTable.getCell(2, 2);\nlocalidadCell.getRange().insertText(\"\", \"Replace\");\n\n// 3) Delete the whole row for \"Jim\u00e9nez Coello Daniel\" \u2014 row index 3 (0 = header).\nconst rosterRows = rosterTable.rows;\nrosterRows.load(\"items\");\nawait context.sync();\nrosterRows.items[3].delete();\n\n// 4) Fix the date in the signature line.\nconst signMatches = body.search(\"a  11  de marzo  2022\", { matchCase: true });\nsignMatches.load(\"items\");\nawait context.sync();\nif (signMatches.items.length > 0) {\n  signMatches.items[0].insertText(\"a  21  de marzo  2022\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop script. $word.ActiveDocument is the open document.\n#\n# Target edit (per the commit's XML diff):\n#  1. \"...suscrito con fecha  11 de marzo de 2022 entre...\"  -> \"...21 de marzo de 2022...\"\n#  2. The table cell holding \"Almod\u00f3var del Campo\" (LOCALIDAD DE RESIDENCIA,\n#     row for \"D\u00edez Vi\u00f1as Malena\") is emptied out.\n#  3. The entire table row for \"Jim\u00e9nez Coello Daniel\" is removed.\n#  4. \"En Puertollano a  11  de marzo  2022\" -> \"En Puertollano a  21  de marzo  2022\"\n\n$d = $word.ActiveDocument\n\n# 1) Fix the date in the intro paragraph (unique match incl. surrounding text).\n$introRange = $d.Content\n$introRange.Find.Execute(\"fecha  11 de marzo de 2022\", $false, $false, $false, $false, $false, $true, 1, $false, \"fecha  21 de marzo de 2022\", 2)\n\n# 2) & 3) Work on the student-roster table (2nd table in the document body).\n$rosterTable = $d.Tables.Item(2)\n\n# 2) Empty the \"Almod\u00f3var del Campo\" cell \u2014 row 3 (\"D\u00edez Vi\u00f1as Malena\"),\n#    column 3 (\"LOCALIDAD DE RESIDENCIA DEL ALUMNO/A (**)\"). (1-indexed)\n$localidadCell = $rosterTable.Cell(3, 3)\n$localidadCell.Range.Text = \"\"\n\n# 3) Delete the whole row for \"Jim\u00e9nez Coello Daniel\" \u2014 row 4 (1 = header row).\n$rosterTable.Rows.Item(4).Delete()\n\n# 4) Fix the date in the signature line.\n$signRange = $d.Content\n$signRange.Find.Execute(\"a  11  de marzo  2022\", $false, $false, $false, $false, $false, $true, 1, $false, \"a  21  de marzo  2022\", 2)\n"}
